$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-23 with new parameter values/formulas ---
    # Row 2
    $ws.Range("A2").Value2 = "PAct1_LacI"
    $ws.Range("B2").Value2 = 0.01
    $ws.Range("C2").Formula = "=B2*1000"
    $ws.Range("D2").Value2 = 1.0
    $ws.Range("E2").Value2 = 98.0
    $ws.Range("F2").Value2 = "yes"
    $ws.Range("G2").Value2 = "k_{LacI}"
    # Row 3
    $ws.Range("A3").Value2 = "P4Lacn_cit"
    $ws.Range("B3").Value2 = 0.01
    $ws.Range("C3").Formula = "=B3*1000"
    $ws.Range("D3").Value2 = 1.0
    $ws.Range("E3").Value2 = 92.0
    $ws.Range("F3").Value2 = "yes"
    $ws.Range("G3").Value2 = "k_{Cit}"
    # Row 4
    $ws.Range("A4").Value2 = "dLacI"
    $ws.Range("B4").Value2 = 0.001
    $ws.Range("C4").Formula = "=B4*1000"
    $ws.Range("D4").Value2 = 1.0
    $ws.Range("E4").Value2 = 0.6
    $ws.Range("F4").Value2 = "yes"
    $ws.Range("G4").Value2 = "d_{LacI}"
    # Row 5
    $ws.Range("A5").Value2 = "dCit"
    $ws.Range("B5").Value2 = 0.001
    $ws.Range("C5").Formula = "=B5*1000"
    $ws.Range("D5").Value2 = 0.0
    $ws.Range("E5").Value2 = 0.0077
    $ws.Range("F5").Value2 = "yes"
    $ws.Range("G5").Value2 = "d_{Cit}"
    # Row 6
    $ws.Range("A6").Value2 = "LacI_rep_WT"
    $ws.Range("B6").Value2 = 0.0001
    $ws.Range("C6").Formula = "=B6*1000"
    $ws.Range("D6").Value2 = 1.0
    $ws.Range("E6").Value2 = 1.0
    $ws.Range("F6").Value2 = "yes"
    $ws.Range("G6").Value2 = "theta_{LacI}"
    # Row 7
    $ws.Range("A7").Value2 = "KdLacI"
    $ws.Range("B7").Value2 = 0.01
    $ws.Range("C7").Formula = "=B7*1000"
    $ws.Range("D7").Value2 = 1.0
    $ws.Range("E7").Value2 = 0.9
    $ws.Range("F7").Value2 = "yes"
    $ws.Range("G7").Value2 = "Kd"
    # Row 8
    $ws.Range("A8").Value2 = "nLacI"
    $ws.Range("B8").Value2 = 1.0
    $ws.Range("C8").Value2 = 10.0
    $ws.Range("D8").Value2 = 0.0
    $ws.Range("E8").Value2 = 1.0
    $ws.Range("F8").Value2 = "yes"
    $ws.Range("G8").Value2 = "n"
    # Row 9
    $ws.Range("A9").Value2 = "mu"
    $ws.Range("B9").Value2 = 0.0
    $ws.Range("C9").Formula = "=B9*1000"
    $ws.Range("D9").Value2 = 0.0
    $ws.Range("E9").Value2 = 0.0077
    $ws.Range("F9").Value2 = "no"
    $ws.Range("G9").Value2 = "mu"
    # Row 10
    $ws.Range("A10").Value2 = "kmaturation"
    $ws.Range("B10").Value2 = 0.0
    $ws.Range("C10").Formula = "=B10*1000"
    $ws.Range("D10").Value2 = 0.0
    $ws.Range("E10").Value2 = 0.0173
    $ws.Range("F10").Value2 = "no"
    $ws.Range("G10").Value2 = "kmaturation"
    # Row 11
    $ws.Range("A11").Value2 = "nMperUnit"
    $ws.Range("B11").Value2 = 0.1
    $ws.Range("C11").Formula = "=B11*1000"
    $ws.Range("D11").Value2 = 1.0
    $ws.Range("E11").Value2 = 1.0
    $ws.Range("F11").Value2 = "yes"
    $ws.Range("G11").Value2 = "nMperUnit"
    # Row 12
    $ws.Range("A12").Value2 = "LacI_rep_W220F"
    $ws.Range("B12").Value2 = 0.001
    $ws.Range("C12").Formula = "=B12*1000"
    $ws.Range("D12").Value2 = 1.0
    $ws.Range("E12").Value2 = 1.0
    $ws.Range("F12").Value2 = "yes"
    $ws.Range("G12").Value2 = "theta_{LacI_W220F}"
    # Row 13
    $ws.Range("A13").Value2 = "IPTG"
    $ws.Range("B13").Value2 = 0.0
    $ws.Range("C13").Formula = "=B13*1000"
    $ws.Range("D13").Value2 = 0.0
    $ws.Range("E13").Value2 = 250.0
    $ws.Range("F13").Value2 = "no"
    $ws.Range("G13").Value2 = "IPTGAdded"
    # Row 14
    $ws.Range("A14").Value2 = "indTime"
    $ws.Range("B14").Value2 = 0.0
    $ws.Range("C14").Formula = "=B14*1000"
    $ws.Range("D14").Value2 = 0.0
    $ws.Range("E14").Value2 = 5000.0
    $ws.Range("F14").Value2 = "no"
    $ws.Range("G14").Value2 = "ind-time"
    # Row 15
    $ws.Range("A15").Value2 = "P_4Lacn_LacI"
    $ws.Range("B15").Value2 = 0.01
    $ws.Range("C15").Formula = "=B15*1000"
    $ws.Range("D15").Value2 = 1.0
    $ws.Range("E15").Value2 = 98.0
    $ws.Range("F15").Value2 = "yes"
    $ws.Range("G15").Value2 = "k_{LacI_W220F_Q60G_T167A}"
    # Row 16
    $ws.Range("A16").Value2 = "P_4Lacn_LacI_L"
    $ws.Range("B16").Value2 = 0.00001
    $ws.Range("C16").Formula = "=B16*1000"
    $ws.Range("D16").Value2 = 1.0
    $ws.Range("E16").Value2 = 0.0003
    $ws.Range("F16").Value2 = "yes"
    $ws.Range("G16").Value2 = "kL_W220F_Q60G_T167A"
    # Row 17
    $ws.Range("A17").Value2 = "LacI_rep_3mut"
    $ws.Range("B17").Value2 = 0.01
    $ws.Range("C17").Formula = "=B17*1000"
    $ws.Range("D17").Value2 = 1.0
    $ws.Range("E17").Value2 = 1.0
    $ws.Range("F17").Value2 = "yes"
    $ws.Range("G17").Value2 = "theta_{LacI_W220F_Q60G_T167A}"
    # Row 18
    $ws.Range("A18").Value2 = "Silence_LacI_rep"
    $ws.Range("B18").Value2 = 0.0
    $ws.Range("C18").Value2 = 1.0
    $ws.Range("D18").Value2 = 0.0
    $ws.Range("E18").Value2 = 1.0
    $ws.Range("F18").Value2 = "no"
    $ws.Range("G18").Value2 = "silence"
    # Row 19
    $ws.Range("A19").Value2 = "pt7_LacI"
    $ws.Range("B19").Value2 = 0.01
    $ws.Range("C19").Formula = "=B19*1000"
    $ws.Range("D19").Value2 = 1.0
    $ws.Range("E19").Value2 = 98.0
    $ws.Range("F19").Value2 = "yes"
    $ws.Range("G19").Value2 = "k_{pt7_PacI}"
    # Row 20
    $ws.Range("A20").Value2 = "P3_Lacn_5_cit"
    $ws.Range("B20").Value2 = 0.01
    $ws.Range("C20").Formula = "=B20*1000"
    $ws.Range("D20").Value2 = 1.0
    $ws.Range("E20").Value2 = 92.0
    $ws.Range("F20").Value2 = "yes"
    $ws.Range("G20").Value2 = "k_{Cit_Lacn3}"
    # Row 21
    $ws.Range("A21").Value2 = "P3_Lacn_5_cit_L"
    $ws.Range("B21").Value2 = 0.00001
    $ws.Range("C21").Formula = "=B21*1000"
    $ws.Range("D21").Value2 = 1.0
    $ws.Range("E21").Value2 = 0.0003
    $ws.Range("F21").Value2 = "yes"
    $ws.Range("G21").Value2 = "kL_cit_pt7"
    # Row 22
    $ws.Range("A22").Value2 = "dLacI_pt7"
    $ws.Range("B22").Value2 = 0.01
    $ws.Range("C22").Formula = "=B22*1000"
    $ws.Range("D22").Value2 = 1.0
    $ws.Range("E22").Value2 = 0.6
    $ws.Range("F22").Value2 = "yes"
    $ws.Range("G22").Value2 = "d_{LacI_pt7}"
    # Row 23
    $ws.Range("A23").Value2 = "nLacI_P3"
    $ws.Range("B23").Value2 = 1.0
    $ws.Range("C23").Value2 = 10.0
    $ws.Range("D23").Value2 = 0.0
    $ws.Range("E23").Value2 = 1.0
    $ws.Range("F23").Value2 = "yes"
    $ws.Range("G23").Value2 = "n_P3"

# --- Remove the now-obsolete last row (row 24) ---
$ws.Rows(24).Delete()

# --- Column C: best-fit width for the new formula column ---
$ws.Columns("C").ColumnWidth = 11.33

# --- Restore the active selection shown in the sheet view ---
$ws.Range("B22").Select()
